$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data fetched from Adafruit IO, identical in shape/content to
# the previous reading (row 67), so copy that row's formatting/types down
# to the new row 68 rather than typing values in (which would coerce the
# numeric-looking "25" into a real number instead of text).
$ws.Range("A67:F67").Copy()
$ws.Range("A68:F68").PasteSpecial(-4104)
$excel.CutCopyMode = $false
